$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AN3").Value = 2.74
$ws.Range("I3").Value = 870
$ws.Range("J3").Value = 9.199999999999999
$ws.Range("K3").Value = 950
$ws.Range("N3").Value = 1.1
$ws.Range("O3").Value = 1.09
$ws.Range("P3").Value = 3.45
$ws.Range("Q3").Value = 1.27
$ws.Range("R3").Value = 2.02
$ws.Range("S3").Value = 1.67
$ws.Range("T3").Value = 2.02
$ws.Range("U3").Value = 1.78
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 12
$ws.Range("AK4").Value = 55
$ws.Range("AL4").Value = 70
$ws.Range("AN4").Value = 65
$ws.Range("I4").Value = 2.38
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 980
$ws.Range("F5").Value = 1.78
$ws.Range("N5").Value = 2.74
$ws.Range("T5").Value = 2.16
$ws.Range("AA6").Value = 50
$ws.Range("G6").Value = 2.8
$ws.Range("I6").Value = 2.92
$ws.Range("J6").Value = 3.45
$ws.Range("T6").Value = 1.44
$ws.Range("V6").Value = 1.52
$ws.Range("G9").Value = 17
$ws.Range("I9").Value = 1.39
$ws.Range("J9").Value = 4.4
$ws.Range("N9").Value = 3.3
$ws.Range("Q9").Value = 1.9
$ws.Range("V9").Value = 3.5
$ws.Range("AL10").Value = 60
$ws.Range("AA11").Value = 13
$ws.Range("AD11").Value = 12.5
$ws.Range("AE11").Value = 14.5
$ws.Range("AG11").Value = 980
$ws.Range("AH11").Value = 26
$ws.Range("AO11").Value = 3.7
$ws.Range("F11").Value = 8.199999999999999
$ws.Range("G11").Value = 11.5
$ws.Range("H11").Value = 1.28
$ws.Range("I11").Value = 1.36
$ws.Range("J11").Value = 5.7
$ws.Range("K11").Value = 7.8
$ws.Range("P11").Value = 3.25
$ws.Range("Q11").Value = 1.35
$ws.Range("R11").Value = 1.92
$ws.Range("T11").Value = 1.7
$ws.Range("V11").Value = 3.55
$ws.Range("W11").Value = 1.1
$ws.Range("Y11").Value = 16.5
$ws.Range("Z11").Value = 12.5
$ws.Range("F12").Value = 2.32
$ws.Range("G12").Value = 2.38
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 3.4
$ws.Range("J12").Value = 3.65
$ws.Range("K12").Value = 3.7
$ws.Range("L12").Value = 1.34
$ws.Range("P12").Value = 2.26
$ws.Range("Q12").Value = 1.74
$ws.Range("R12").Value = 1.51
$ws.Range("S12").Value = 2.76
$ws.Range("T12").Value = 1.61
$ws.Range("U12").Value = 2.46
$ws.Range("V12").Value = 1.41
$ws.Range("W12").Value = 1.73
$ws.Range("AE14").Value = 85
$ws.Range("F14").Value = 1.86
$ws.Range("G14").Value = 1.95
$ws.Range("H14").Value = 4.5
$ws.Range("I14").Value = 5.1
$ws.Range("J14").Value = 3.65
$ws.Range("V14").Value = 1.25
$ws.Range("W14").Value = 2.04
$ws.Range("N15").Value = 7
$ws.Range("P15").Value = 3.05
$ws.Range("R15").Value = 1.84
$ws.Range("S15").Value = 2.06
$ws.Range("T15").Value = 1.79
$ws.Range("W15").Value = 3.95
$ws.Range("AB16").Value = 980
$ws.Range("F16").Value = 3.35
$ws.Range("J16").Value = 3.55
$ws.Range("K16").Value = 3.9
$ws.Range("L16").Value = 1.31
$ws.Range("O16").Value = 1.24
$ws.Range("T16").Value = 1.68
$ws.Range("U16").Value = 2.2
$ws.Range("W16").Value = 1.35
$ws.Range("AA17").Value = 38
$ws.Range("AB17").Value = 22
$ws.Range("AC17").Value = 12
$ws.Range("AD17").Value = 14.5
$ws.Range("AE17").Value = 27
$ws.Range("AF17").Value = 30
$ws.Range("AG17").Value = 17
$ws.Range("AH17").Value = 18
$ws.Range("AI17").Value = 36
$ws.Range("AJ17").Value = 60
$ws.Range("AK17").Value = 36
$ws.Range("AL17").Value = 40
$ws.Range("AM17").Value = 70
$ws.Range("AN17").Value = 23
$ws.Range("AO17").Value = 15
$ws.Range("N17").Value = 5.3
$ws.Range("O17").Value = 1.18
$ws.Range("T17").Value = 1.54
$ws.Range("U17").Value = 2.52
$ws.Range("W17").Value = 1.43
$ws.Range("X17").Value = 30
$ws.Range("Y17").Value = 18.5
$ws.Range("Z17").Value = 23
$ws.Range("AC18").Value = 15.5
$ws.Range("AD18").Value = 15
$ws.Range("AF18").Value = 42
$ws.Range("AH18").Value = 17.5
$ws.Range("AJ18").Value = 70
$ws.Range("AM18").Value = 48
$ws.Range("AN18").Value = 18
$ws.Range("AO18").Value = 8.199999999999999
$ws.Range("F18").Value = 3.25
$ws.Range("I18").Value = 2.22
$ws.Range("J18").Value = 3.8
$ws.Range("N18").Value = 6.8
$ws.Range("O18").Value = 1.11
$ws.Range("R18").Value = 1.96
$ws.Range("S18").Value = 1.87
$ws.Range("T18").Value = 1.39
$ws.Range("U18").Value = 3
$ws.Range("V18").Value = 1.84
$ws.Range("X18").Value = 55
$ws.Range("AA19").Value = 44
$ws.Range("AB19").Value = 12.5
$ws.Range("AC19").Value = 8.6
$ws.Range("AD19").Value = 15
$ws.Range("AE19").Value = 42
$ws.Range("AF19").Value = 32
$ws.Range("AG19").Value = 20
$ws.Range("AH19").Value = 28
$ws.Range("AI19").Value = 75
$ws.Range("AJ19").Value = 110
$ws.Range("AK19").Value = 75
$ws.Range("AL19").Value = 100
$ws.Range("AM19").Value = 220
$ws.Range("AN19").Value = 110
$ws.Range("AO19").Value = 44
$ws.Range("F19").Value = 3.4
$ws.Range("G19").Value = 4.3
$ws.Range("H19").Value = 2.24
$ws.Range("I19").Value = 2.6
$ws.Range("J19").Value = 2.84
$ws.Range("K19").Value = 3.45
$ws.Range("L19").Value = 1.47
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 2.52
$ws.Range("O19").Value = 1.49
$ws.Range("P19").Value = 1.51
$ws.Range("Q19").Value = 2.46
$ws.Range("S19").Value = 4.5
$ws.Range("T19").Value = 2.02
$ws.Range("U19").Value = 1.75
$ws.Range("V19").Value = 1.62
$ws.Range("W19").Value = 1.31
$ws.Range("X19").Value = 11
$ws.Range("Y19").Value = 9
$ws.Range("Z19").Value = 17
$ws.Range("F20").Value = 1.98
$ws.Range("G20").Value = 2.56
$ws.Range("I20").Value = 4.9
$ws.Range("N20").Value = 1.01
$ws.Range("O20").Value = 1.41
$ws.Range("P20").Value = 1.24
$ws.Range("Q20").Value = 1.41
$ws.Range("S20").Value = 1.01
$ws.Range("V20").Value = 1.25
$ws.Range("W20").Value = 1.64
$ws.Range("H21").Value = 9
$ws.Range("I21").Value = 9.199999999999999
$ws.Range("L21").Value = 1.44
$ws.Range("W21").Value = 2.9
$ws.Range("AH23").Value = 20
$ws.Range("H23").Value = 4.3
$ws.Range("P23").Value = 1.79
$ws.Range("T23").Value = 1.98
$ws.Range("U23").Value = 1.98
$ws.Range("W23").Value = 1.92
$ws.Range("X23").Value = 11.5
$ws.Range("G24").Value = 7.4
$ws.Range("H24").Value = 1.57
$ws.Range("I24").Value = 1.59
$ws.Range("J24").Value = 4.3
$ws.Range("L24").Value = 1.36
$ws.Range("N24").Value = 4.1
$ws.Range("O24").Value = 1.3
$ws.Range("Q24").Value = 1.88
$ws.Range("S24").Value = 3.25
$ws.Range("T24").Value = 1.97
$ws.Range("V24").Value = 2.68
$ws.Range("W24").Value = 1.16
$ws.Range("F25").Value = 1.88
$ws.Range("N25").Value = 3.9
$ws.Range("AA26").Value = 270
$ws.Range("AB26").Value = 7.6
$ws.Range("AC26").Value = 10
$ws.Range("AE26").Value = 160
$ws.Range("AF26").Value = 11.5
$ws.Range("AG26").Value = 13
$ws.Range("AI26").Value = 170
$ws.Range("AJ26").Value = 24
$ws.Range("AK26").Value = 29
$ws.Range("AL26").Value = 70
$ws.Range("AM26").Value = 280
$ws.Range("AN26").Value = 22
$ws.Range("AO26").Value = 270
$ws.Range("I26").Value = 6.8
$ws.Range("L26").Value = 1.53
$ws.Range("M26").Value = 1.09
$ws.Range("N26").Value = 2.74
$ws.Range("O26").Value = 1.49
$ws.Range("P26").Value = 1.59
$ws.Range("R26").Value = 1.21
$ws.Range("S26").Value = 4.8
$ws.Range("T26").Value = 2.24
$ws.Range("U26").Value = 1.7
$ws.Range("W26").Value = 2.18
$ws.Range("X26").Value = 12
$ws.Range("Y26").Value = 19
$ws.Range("Z26").Value = 60
$ws.Range("G27").Value = 1.74
$ws.Range("N27").Value = 2.62
$ws.Range("Q27").Value = 1.91
$ws.Range("S27").Value = 3.05
$ws.Range("T27").Value = 1.01
$ws.Range("W27").Value = 2.34
